# Generate Report for Handoff
#
# A new handoff/XLIFF-generate cycle ran for the last tracked file
# (da9e7863-03e8-4548-a100-4732be3f6675), so its recorded timestamps move
# forward on every sheet that tracks it:
#   - Overview!G7            "Latest HO Xliff Generate Date"   -> 2016-09-08 04:55:53
#   - zh-cn!H7 (row 7)       "Latest Handoff Datetime"         -> 2016-09-08 04:55:48
#   - de-de!H7 (row 7)       "Latest Handoff Datetime"         -> 2016-09-08 04:55:53

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-09-08 04:55:53"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-09-08 04:55:48"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-09-08 04:55:53"
